# "Colocando header nos graficos" - add a header label to column A on each
# data sheet, fix missing accents in the Portuguese labels, and drop the
# unused "Teto" row / rename + renumber the "Custo Total" sheet.

$wb = $excel.ActiveWorkbook

$xlPasteAll     = -4104
$xlPasteFormats = -4122

function Set-ColumnAHeader {
    param($ws, [string]$text)
    # New header cell in A1. Copy the formatting (bold/border/centered)
    # that already lives on B1 so the new cell matches the rest of the
    # header row without minting a brand-new style entry.
    $ws.Range("A1").Value = $text
    $ws.Range("B1").Copy()
    $ws.Range("A1").PasteSpecial($xlPasteFormats)
}

function Set-LabelNoStyle {
    param($ws, [string]$cellAddr, $newText)
    # Row labels in column A lose their bold/border styling; reuse the
    # (unstyled) numeric cell B2's formatting to clear it cleanly.
    if ($null -ne $newText) {
        $ws.Range($cellAddr).Value = $newText
    }
    $ws.Range("B2").Copy()
    $ws.Range($cellAddr).PasteSpecial($xlPasteFormats)
}

# --- Sheets 1-4: same Fonte/Tecnologia table layout ------------------------
for ($idx = 1; $idx -le 4; $idx++) {
    $ws = $wb.Worksheets.Item($idx)

    Set-ColumnAHeader $ws "Fonte/Tecnologia"

    Set-LabelNoStyle $ws "A2"  $null
    Set-LabelNoStyle $ws "A3"  "Gás Natural"
    Set-LabelNoStyle $ws "A4"  "Carvão"
    Set-LabelNoStyle $ws "A5"  $null
    Set-LabelNoStyle $ws "A6"  "Óleos Comb"
    Set-LabelNoStyle $ws "A7"  $null
    Set-LabelNoStyle $ws "A8"  "Eólica"
    Set-LabelNoStyle $ws "A9"  $null
    Set-LabelNoStyle $ws "A10" $null
    Set-LabelNoStyle $ws "A11" "Pot. Compl."
    Set-LabelNoStyle $ws "A12" $null
}

# --- Sheet 5: Emissoes Totais ----------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

Set-ColumnAHeader $ws5 "Período"
Set-LabelNoStyle $ws5 "A2" "P.Médio"
Set-LabelNoStyle $ws5 "A3" "P.Crítico"
$ws5.Rows.Item(4).Delete()

# --- Sheet 6: Custo Total ---------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)
$ws1 = $wb.Worksheets.Item(1)

$ws6.Range("A1").Value = "Tipo Expansão"
$ws6.Range("B1").Copy()
$ws6.Range("A1").PasteSpecial($xlPasteFormats)

# Reuse sheet1's "2015" header cell (text + style) verbatim for B1 so it
# stays a real text value instead of Excel auto-coercing "2015" to a number.
$ws1.Range("B1").Copy()
$ws6.Range("B1").PasteSpecial($xlPasteAll)

Set-LabelNoStyle $ws6 "A2" "Expansão Centralizada"
$ws6.Range("B2").Value = 604

Set-LabelNoStyle $ws6 "A3" "Expansão por GD"
$ws6.Range("B3").Value = 99

Write-Output "Edits applied"
